$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 5-19 (A: archivo, B: nombre_detectado)
$data = @(
    @("0479f849-3354-4b13-bc69-7b30043bd56e_fila_10.png", "Abraham Silva Ampre"),
    @("315628bb-db42-48bc-8fe3-31ea92794346_fila_1.png", "Hotep Antonio Ruiz Lezama"),
    @("3a478992-9152-48ad-ae2c-d8a0fc05e7ba_fila_13.png", "Marlon Josue Gonzales Cano"),
    @("3c943684-5c12-4960-b9aa-75089d79a55d_fila_4.png", "Yadder Fernando Torres"),
    @("3f3a7c2f-0abc-47fc-af0b-accf90aed4d5_fila_9.png", "Roman Alfonso Grios Boza"),
    @("57df3fb7-62cb-4e4d-bce7-5ef755df64f1_fila_6.png", "Ronier Jose Rivera"),
    @("648ad02a-46de-4ec7-b602-adc46e8705eb_fila_2.png", "Isabella Dompe Estrada"),
    @("676ad72f-8b15-4c34-b497-3fcc0fdc1249_fila_7.png", "Cristina Jozabed Carvajal"),
    @("89fdf458-71de-4059-b799-073878fac609_fila_14.png", "Angel Isaac Alvarez Quiñonez"),
    @("b41224c4-a34d-48fa-a1cc-287876b653ea_fila_8.png", "David Orlando Mena Valverde"),
    @("c8324d8f-0bbc-458f-b62a-6631c96ce433_fila_5.png", "Erick Espinoza"),
    @("de39c47b-8c5a-4f19-b6bf-508419360394_fila_3.png", "Bryan Alexander Cano"),
    @("e774e6cd-a59a-41dc-b8e5-754879731161_fila_12.png", "José Danilo Suárez"),
    @("f16b2728-ccac-4f7a-84c0-c1333af701e8_fila_11.png", "Eduardo Domingo Zeledon Merca"),
    @("fe13f750-c24d-4cd3-b041-3d76d19e3c54_fila_15.png", "Orlando Mauricio Guevara")
)

$startRow = 5
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
